# Fix Table & Loc, can run in android
#
# 1) Row 11: the two image-path strings change content
#    (icon.png -> s2.png, btn_press.png -> s1.png) and the
#    B11/C11 cells end up pointing at the (now renumbered) shared
#    strings 27/28 respectively.
# 2) A brand new row 12 is added with a new LOCIMAGE_VAL_TEXTUREA
#    entry, pointing a Texture path (used for both B12 and C12).
# 3) The sheet's used range / dimension grows to A1:S12 and the
#    active selection moves to C13.
# 4) Conditional formatting (the same "formula highlight" rule
#    already used on the sheet) is extended to cover the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) update the two existing image-path strings in row 11 ---
$ws.Cells.Item(11, 2).Value = "Assets/Res/UI/Sprite/s2.png"
$ws.Cells.Item(11, 3).Value = "Assets/Res/UI/Sprite/s1.png"

# --- 2) add row 12, copying the look (style) of row 11 first ---
$ws.Range("A12:C12").Style = $ws.Range("A11").Style

$ws.Cells.Item(12, 1).Value = "LOCIMAGE_VAL_TEXTUREA"
$ws.Cells.Item(12, 2).Value = "Assets/Res/UI/Texture/t2.png"
$ws.Cells.Item(12, 3).Value = "Assets/Res/UI/Texture/t2.png"

# --- 3) move the selection the way it ends up after the edit ---
$ws.Range("C13").Select()

# --- 4) extend the existing conditional formatting ("highlight
#        formula cells") rule to also cover the new row (B12 and
#        C12), reusing the same fill colour as the pre-existing
#        rule. ---
$cf = $ws.Range("B12:C12").FormatConditions.Add(2, 0, "公式单元格=TRUE")
$cf.Interior.Color = 15917529
